{"js": "// Apply the text replacements described by the diff:\n// - the date line at the top of the document\n// - 25 \"AxB=\" multiplication prompts scattered through the table cells\n//\n// Every old string in the mapping is unique within the document, so a\n// plain body.search(...) + insertText(..., \"Replace\") per pair is safe\n// and unambiguous.\n\nconst replacements = [\n  [\"2024-08-21 Wednesday\", \"2024-08-22 Thursday\"],\n  [\"728\u00d76=\", \"572\u00d78=\"],\n  [\"897\u00d75=\", \"346\u00d77=\"],\n  [\"429\u00d75=\", \"571\u00d73=\"],\n  [\"229\u00d73=\", \"108\u00d75=\"],\n  [\"748\u00d76=\", \"594\u00d74=\"],\n  [\"445\u00d78=\", \"665\u00d77=\"],\n  [\"157\u00d77=\", \"955\u00d78=\"],\n  [\"307\u00d74=\", \"863\u00d72=\"],\n  [\"990\u00d76=\", \"575\u00d73=\"],\n  [\"647\u00d77=\", \"449\u00d76=\"],\n  [\"320\u00d78=\", \"701\u00d75=\"],\n  [\"605\u00d76=\", \"908\u00d73=\"],\n  [\"169\u00d77=\", \"578\u00d79=\"],\n  [\"531\u00d75=\", \"389\u00d72=\"],\n  [\"696\u00d75=\", \"769\u00d73=\"],\n  [\"436\u00d74=\", \"534\u00d72=\"],\n  [\"424\u00d78=\", \"398\u00d76=\"],\n  [\"721\u00d77=\", \"310\u00d79=\"],\n  [\"807\u00d72=\", \"799\u00d78=\"],\n  [\"344\u00d76=\", \"858\u00d76=\"],\n  [\"832\u00d74=\", \"849\u00d72=\"],\n  [\"635\u00d79=\", \"190\u00d77=\"],\n  [\"322\u00d73=\", \"767\u00d75=\"],\n  [\"247\u00d73=\", \"159\u00d77=\"],\n  [\"278\u00d79=\", \"271\u00d74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Search text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the text replacements described by the diff:\n# - the date line at the top of the document\n# - 25 \"AxB=\" multiplication prompts scattered through the table cells\n#\n# Every old string in the mapping is unique within the document, so a\n# plain Find/Replace (wdReplaceAll, scoped to the whole story) per pair\n# is safe and unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-08-21 Wednesday\", \"2024-08-22 Thursday\"),\n    @(\"728\u00d76=\", \"572\u00d78=\"),\n    @(\"897\u00d75=\", \"346\u00d77=\"),\n    @(\"429\u00d75=\", \"571\u00d73=\"),\n    @(\"229\u00d73=\", \"108\u00d75=\"),\n    @(\"748\u00d76=\", \"594\u00d74=\"),\n    @(\"445\u00d78=\", \"665\u00d77=\"),\n    @(\"157\u00d77=\", \"955\u00d78=\"),\n    @(\"307\u00d74=\", \"863\u00d72=\"),\n    @(\"990\u00d76=\", \"575\u00d73=\"),\n    @(\"647\u00d77=\", \"449\u00d76=\"),\n    @(\"320\u00d78=\", \"701\u00d75=\"),\n    @(\"605\u00d76=\", \"908\u00d73=\"),\n    @(\"169\u00d77=\", \"578\u00d79=\"),\n    @(\"531\u00d75=\", \"389\u00d72=\"),\n    @(\"696\u00d75=\", \"769\u00d73=\"),\n    @(\"436\u00d74=\", \"534\u00d72=\"),\n    @(\"424\u00d78=\", \"398\u00d76=\"),\n    @(\"721\u00d77=\", \"310\u00d79=\"),\n    @(\"807\u00d72=\", \"799\u00d78=\"),\n    @(\"344\u00d76=\", \"858\u00d76=\"),\n    @(\"832\u00d74=\", \"849\u00d72=\"),\n    @(\"635\u00d79=\", \"190\u00d77=\"),\n    @(\"322\u00d73=\", \"767\u00d75=\"),\n    @(\"247\u00d73=\", \"159\u00d77=\"),\n    @(\"278\u00d79=\", \"271\u00d74=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($old, $true, $false, $false, $false, $false, $true, 0, $false, $new, 2) | Out-Null\n}\n"}
